$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 17.04695337284734
$ws.Range("D2").Value = 9.437467674464822
$ws.Range("E2").Value = 15.39472393693095
$ws.Range("F2").Value = 36.14610799108036
$ws.Range("G2").Value = 3.674760059667096
$ws.Range("J2").Value = 11.39301999968659
$ws.Range("K2").Value = 10.19176106570954
$ws.Range("L2").Value = 9.237302152776941
$ws.Range("N2").Value = 21.13380495858535
$ws.Range("O2").Value = 27.17288645156798
$ws.Range("B3").Value = 16.9528335583414
$ws.Range("D3").Value = 9.443665523333333
$ws.Range("E3").Value = 15.43063691350818
$ws.Range("F3").Value = 36.21458280853332
$ws.Range("G3").Value = 3.676629297140562
$ws.Range("J3").Value = 11.41533253785539
$ws.Range("K3").Value = 9.976267810906858
$ws.Range("L3").Value = 9.202724197122031
$ws.Range("N3").Value = 21.19588656612919
$ws.Range("O3").Value = 27.23825118037559
$ws.Range("B4").Value = 16.89787549928091
$ws.Range("D4").Value = 9.448594303626768
$ws.Range("E4").Value = 15.45425446794267
$ws.Range("F4").Value = 36.26379712474428
$ws.Range("G4").Value = 3.677839038701456
$ws.Range("J4").Value = 11.42978225410541
$ws.Range("K4").Value = 9.843136988507098
$ws.Range("L4").Value = 9.182625483786726
$ws.Range("N4").Value = 21.23580761379689
$ws.Range("O4").Value = 27.28341373867755
$ws.Range("B5").Value = 16.87621050910846
$ws.Range("D5").Value = 9.4508859451338
$ws.Range("E5").Value = 15.46427349496342
$ws.Range("F5").Value = 36.28565355641388
$ws.Range("G5").Value = 3.678347662376914
$ws.Range("J5").Value = 11.43585966744405
$ws.Range("K5").Value = 9.78875522744829
$ws.Range("L5").Value = 9.174725380468596
$ws.Range("N5").Value = 21.25253039550657
$ws.Range("O5").Value = 27.30308086938874
$ws.Range("B6").Value = 16.8726577247456
$ws.Range("D6").Value = 9.451283588980633
$ws.Range("E6").Value = 15.46596100480155
$ws.Range("F6").Value = 36.2893915261543
$ws.Range("G6").Value = 3.678433065238796
$ws.Range("J6").Value = 11.43688024971309
$ws.Range("K6").Value = 9.77971958697105
$ws.Range("L6").Value = 9.173431263780513
$ws.Range("N6").Value = 21.25533469837555
$ws.Range("O6").Value = 27.30642282130482
$ws.Range("B7").Value = 16.89758033437188
$ws.Range("D7").Value = 9.448624062294151
$ws.Range("E7").Value = 15.4543879890702
$ws.Range("F7").Value = 36.26408459800617
$ws.Range("G7").Value = 3.677845834771223
$ws.Range("J7").Value = 11.42986345014152
$ws.Range("K7").Value = 9.842404001351243
$ws.Range("L7").Value = 9.1825177577403
$ws.Range("N7").Value = 21.23603130053503
$ws.Range("O7").Value = 27.28367386451597
$ws.Range("B8").Value = 17.01392466161737
$ws.Range("D8").Value = 9.439371969526759
$ws.Range("E8").Value = 15.40678194662973
$ws.Range("F8").Value = 36.16822884015064
$ws.Range("G8").Value = 3.675391729084486
$ws.Range("J8").Value = 11.40055806534373
$ws.Range("K8").Value = 10.11767319082579
$ws.Range("L8").Value = 9.225148085565102
$ws.Range("N8").Value = 21.15483735253488
$ws.Range("O8").Value = 27.19437977459022
$ws.Range("B9").Value = 17.2636760782134
$ws.Range("D9").Value = 9.430110886655719
$ws.Range("E9").Value = 15.32582802067256
$ws.Range("F9").Value = 36.03722277636246
$ws.Range("G9").Value = 3.6710691570647
$ws.Range("J9").Value = 11.34901567047783
$ws.Range("K9").Value = 10.64780129922915
$ws.Range("L9").Value = 9.317474504105268
$ws.Range("N9").Value = 21.00985814488675
$ws.Range("O9").Value = 27.05923040850213
$ws.Range("B10").Value = 17.45911335772434
$ws.Range("D10").Value = 9.428678974133874
$ws.Range("E10").Value = 15.27386901990352
$ws.Range("F10").Value = 35.97578135055628
$ws.Range("G10").Value = 3.668188961288628
$ws.Range("J10").Value = 11.31472697259383
$ws.Range("K10").Value = 11.02733312550062
$ws.Range("L10").Value = 9.390269494020769
$ws.Range("N10").Value = 20.91193645222554
$ws.Range("O10").Value = 26.98436792143038
$ws.Range("B11").Value = 17.55034023761887
$ws.Range("D11").Value = 9.42918398902699
$ws.Range("E11").Value = 15.25185500923744
$ws.Range("F11").Value = 35.95539754960836
$ws.Range("G11").Value = 3.666942219228384
$ws.Range("J11").Value = 11.29989845470368
$ws.Range("K11").Value = 11.19704288218792
$ws.Range("L11").Value = 9.424382144694833
$ws.Range("N11").Value = 20.8692368499825
$ws.Range("O11").Value = 26.9556288204852
$ws.Range("B12").Value = 17.58519509987576
$ws.Range("D12").Value = 9.429540554552652
$ws.Range("E12").Value = 15.24375149261192
$ws.Range("F12").Value = 35.94876669606217
$ws.Range("G12").Value = 3.666479188414355
$ws.Range("J12").Value = 11.29439342404699
$ws.Range("K12").Value = 11.26082420556669
$ws.Range("L12").Value = 9.437435981743949
$ws.Range("N12").Value = 20.85333161609925
$ws.Range("O12").Value = 26.94551127621183
$ws.Range("B13").Value = 17.57767512555325
$ws.Range("D13").Value = 9.429456424725471
$ws.Range("E13").Value = 15.24548638931996
$ws.Range("F13").Value = 35.95014638201599
$ws.Range("G13").Value = 3.666578507098433
$ws.Range("J13").Value = 11.29557413664267
$ws.Range("K13").Value = 11.24711021581189
$ws.Range("L13").Value = 9.434618669914109
$ws.Range("N13").Value = 20.85674536231507
$ws.Range("O13").Value = 26.94765622040898
$ws.Range("B14").Value = 17.55320171341006
$ws.Range("D14").Value = 9.42921001618965
$ws.Range("E14").Value = 15.2511836675142
$ws.Range("F14").Value = 35.95483022058264
$ws.Range("G14").Value = 3.666903943611952
$ws.Range("J14").Value = 11.29944334640712
$ws.Range("K14").Value = 11.20230024733935
$ws.Range("L14").Value = 9.425453406636571
$ws.Range("N14").Value = 20.86792302909895
$ws.Range("O14").Value = 26.95478110140651
$ws.Range("B15").Value = 17.53825057351779
$ws.Range("D15").Value = 9.429080584844584
$ws.Range("E15").Value = 15.25470370344604
$ws.Range("F15").Value = 35.95784089638659
$ws.Range("G15").Value = 3.667104464492093
$ws.Range("J15").Value = 11.30182768783611
$ws.Range("K15").Value = 11.1747880618362
$ws.Range("L15").Value = 9.419856922633423
$ws.Range("N15").Value = 20.87480404129388
$ws.Range("O15").Value = 26.95924498766871
$ws.Range("B16").Value = 17.45319608865576
$ws.Range("D16").Value = 9.428669159041721
$ws.Range("E16").Value = 15.27534028377164
$ws.Range("F16").Value = 35.97726574506054
$ws.Range("G16").Value = 3.668271712287578
$ws.Range("J16").Value = 11.31571149784522
$ws.Range("K16").Value = 11.01617769003554
$ws.Range("L16").Value = 9.388059605422868
$ws.Range("N16").Value = 20.91476400457038
$ws.Range("O16").Value = 26.98635314328483
$ws.Range("B17").Value = 17.4015951574746
$ws.Range("D17").Value = 9.428712278148804
$ws.Range("E17").Value = 15.28841527335832
$ws.Range("F17").Value = 35.99112032046609
$ws.Range("G17").Value = 3.669004006608234
$ws.Range("J17").Value = 11.32442554700914
$ws.Range("K17").Value = 10.91807726287014
$ws.Range("L17").Value = 9.368803369390214
$ws.Range("N17").Value = 20.93974996069478
$ws.Range("O17").Value = 27.00434540367684
$ws.Range("B18").Value = 17.37213568226795
$ws.Range("D18").Value = 9.42884598600706
$ws.Range("E18").Value = 15.2960884043482
$ws.Range("F18").Value = 35.9998012999019
$ws.Range("G18").Value = 3.669431179963597
$ws.Range("J18").Value = 11.32951010147138
$ws.Range("K18").Value = 10.86137905413709
$ws.Range("L18").Value = 9.357822067584584
$ws.Range("N18").Value = 20.95429499669902
$ws.Range("O18").Value = 27.01519442370102
$ws.Range("B19").Value = 17.36219973961493
$ws.Range("D19").Value = 9.428909990409629
$ws.Range("E19").Value = 15.29871264634866
$ws.Range("F19").Value = 36.00286284180068
$ws.Range("G19").Value = 3.669576841431813
$ws.Range("J19").Value = 11.33124410400878
$ws.Range("K19").Value = 10.84213698153153
$ws.Range("L19").Value = 9.354120418095532
$ws.Range("N19").Value = 20.95924957865148
$ws.Range("O19").Value = 27.01895362133125
$ws.Range("B20").Value = 17.40706557058477
$ws.Range("D20").Value = 9.428696423384174
$ws.Range("E20").Value = 15.28700761436602
$ws.Range("F20").Value = 35.98957176610936
$ws.Range("G20").Value = 3.668925434347483
$ws.Range("J20").Value = 11.32349042519409
$ws.Range("K20").Value = 10.9285489764471
$ws.Range("L20").Value = 9.370843513841413
$ws.Range("N20").Value = 20.93707218666134
$ws.Range("O20").Value = 27.00237830586892
$ws.Range("B21").Value = 17.56038195396931
$ws.Range("D21").Value = 9.42927791335744
$ws.Range("E21").Value = 15.24950392774591
$ws.Range("F21").Value = 35.95342493628635
$ws.Range("G21").Value = 3.666808108849168
$ws.Range("J21").Value = 11.29830387834157
$ws.Range("K21").Value = 11.21547562597791
$ws.Range("L21").Value = 9.428141832143927
$ws.Range("N21").Value = 20.86463271716137
$ws.Range("O21").Value = 26.952667574373
$ws.Range("B22").Value = 17.66237374233011
$ws.Range("D22").Value = 9.430621118781014
$ws.Range("E22").Value = 15.22634925115707
$ws.Range("F22").Value = 35.936142569597
$ws.Range("G22").Value = 3.665477238749727
$ws.Range("J22").Value = 11.28248517692655
$ws.Range("K22").Value = 11.40015173816912
$ws.Range("L22").Value = 9.466379440491878
$ws.Range("N22").Value = 20.8188286249585
$ws.Range("O22").Value = 26.92463993066402
$ws.Range("B23").Value = 17.60778320224231
$ws.Range("D23").Value = 9.429816427130811
$ws.Range("E23").Value = 15.23858344301901
$ws.Range("F23").Value = 35.9447863369363
$ws.Range("G23").Value = 3.666182721032427
$ws.Range("J23").Value = 11.29086930771553
$ws.Range("K23").Value = 11.30186626427439
$ws.Range("L23").Value = 9.445901521063945
$ws.Range("N23").Value = 20.8431346785629
$ws.Range("O23").Value = 26.93919037565772
$ws.Range("B24").Value = 17.40459175145058
$ws.Range("D24").Value = 9.428703251997794
$ws.Range("E24").Value = 15.28764353058694
$ws.Range("F24").Value = 35.99026963787084
$ws.Range("G24").Value = 3.66896093765164
$ws.Range("J24").Value = 11.32391296095256
$ws.Range("K24").Value = 10.92381564125224
$ws.Range("L24").Value = 9.369920885616429
$ws.Range("N24").Value = 20.9382822466484
$ws.Range("O24").Value = 27.00326605770694
$ws.Range("B25").Value = 17.19392228906021
$ws.Range("D25").Value = 9.431669539896188
$ws.Range("E25").Value = 15.34640494574998
$ws.Range("F25").Value = 36.0665547973481
$ws.Range("G25").Value = 3.672186396715794
$ws.Range("J25").Value = 11.36232830566291
$ws.Range("K25").Value = 10.50585536624161
$ws.Range("L25").Value = 9.291599303529228
$ws.Range("N25").Value = 21.04756327844712
$ws.Range("O25").Value = 27.09150609901124
